$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the price cells that would otherwise be auto-parsed as numbers
# (e.g. "325.84") so they are stored as text, matching the source data which
# treats the whole Price column as display strings.
$numberFormatCells = @("D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D18","D19","D20","D22","D23","D24","D26","D28","D29","D30","D31","D32","D34","D35","D36","D37","D38","D39","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $numberFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = '29.501.06'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '1.909.78'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("E4").Value = '  +0.49%  '
$ws.Range("D5").Value = '325.84'
$ws.Range("E5").Value = '  -0.65%  '
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("D7").Value = '0.4844'
$ws.Range("E7").Value = '  +1.24%  '
$ws.Range("D8").Value = '0.4072'
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("D9").Value = '0.08146'
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("D10").Value = '1.012'
$ws.Range("E10").Value = '  +0.02%  '
$ws.Range("D11").Value = '23.47'
$ws.Range("E11").Value = '  +4.67%  '
$ws.Range("D12").Value = '1.905.52'
$ws.Range("E12").Value = '  +1.16%  '
$ws.Range("D13").Value = '6.024'
$ws.Range("E13").Value = '  +1.05%  '
$ws.Range("D14").Value = '7.102'
$ws.Range("E14").Value = '  -1.00%  '
$ws.Range("D15").Value = '90.36'
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("D16").Value = '0.06799'
$ws.Range("E16").Value = '  +3.05%  '
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").Value = '0.00001040'
$ws.Range("E18").Value = '  +0.78%  '
$ws.Range("D19").Value = '17.69'
$ws.Range("E19").Value = '  -0.34%  '
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("D21").Value = '29.520.63'
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").Value = '5.605'
$ws.Range("E22").Value = '  +0.95%  '
$ws.Range("D23").Value = '11.80'
$ws.Range("E23").Value = '  +2.18%  '
$ws.Range("D24").Value = '2.165'
$ws.Range("E24").Value = '  -1.79%  '
$ws.Range("D25").Value = '2.129.30'
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("D26").Value = '154.64'
$ws.Range("E26").Value = '  +0.56%  '
$ws.Range("E27").Value = '  +1.13%  '
$ws.Range("D28").Value = '6.275'
$ws.Range("E28").Value = '  +8.38%  '
$ws.Range("D29").Value = '2.104'
$ws.Range("E29").Value = '  -1.62%  '
$ws.Range("D30").Value = '119.87'
$ws.Range("E30").Value = '  +1.89%  '
$ws.Range("D31").Value = '1.032'
$ws.Range("E31").Value = '  -3.26%  '
$ws.Range("D32").Value = '0.09560'
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("E33").Value = '  +2.67%  '
$ws.Range("D34").Value = '1.395'
$ws.Range("E34").Value = '  -2.08%  '
$ws.Range("D35").Value = '3.549'
$ws.Range("E35").Value = '  -0.66%  '
$ws.Range("D36").Value = '0.02266'
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("D37").Value = '0.06111'
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("D38").Value = '1.170'
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("D39").Value = '0.5949'
$ws.Range("E39").Value = '  +0.99%  '
$ws.Range("D41").Value = '7.911'
$ws.Range("E41").Value = '  -5.37%  '
$ws.Range("D42").Value = '0.1854'
$ws.Range("E42").Value = '  +0.63%  '
$ws.Range("D43").Value = '2.454'
$ws.Range("E43").Value = '  +1.68%  '
$ws.Range("D44").Value = '1.285'
$ws.Range("E44").Value = '  -1.10%  '
$ws.Range("D45").Value = '0.07716'
$ws.Range("E45").Value = '  -0.60%  '
$ws.Range("D46").Value = '12.41'
$ws.Range("E46").Value = '  +2.36%  '
$ws.Range("D47").Value = '0.5571'
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("D48").Value = '1.953'
$ws.Range("E48").Value = '  +1.12%  '
$ws.Range("D49").Value = '114.82'
$ws.Range("E49").Value = '  +1.20%  '
$ws.Range("D50").Value = '72.65'
$ws.Range("E50").Value = '  +1.39%  '
$ws.Range("D51").Value = '1.052'
$ws.Range("E51").Value = '  +1.87%  '
